$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "28.270.82"
$ws.Range("E2").Value = "  +3.60%  "
Set-TextValue $ws.Range("D3") "1.815.55"
$ws.Range("E3").Value = "  +3.82%  "
Set-TextValue $ws.Range("D4") "0.9990"
$ws.Range("E4").Value = "  -1.89%  "
Set-TextValue $ws.Range("D5") "329.58"
$ws.Range("E5").Value = "  +1.79%  "
Set-TextValue $ws.Range("D6") "0.9975"
$ws.Range("E6").Value = "  -1.54%  "
Set-TextValue $ws.Range("D7") "0.4437"
$ws.Range("E7").Value = "  +5.77%  "
Set-TextValue $ws.Range("D8") "0.3733"
$ws.Range("E8").Value = "  +4.34%  "
Set-TextValue $ws.Range("D9") "44.78"
$ws.Range("E9").Value = "  -0.20%  "
Set-TextValue $ws.Range("D10") "0.07703"
$ws.Range("E10").Value = "  +5.34%  "
Set-TextValue $ws.Range("D11") "1.125"
$ws.Range("E11").Value = "  +1.54%  "
Set-TextValue $ws.Range("D12") "0.9979"
$ws.Range("E12").Value = "  -1.73%  "
Set-TextValue $ws.Range("D13") "22.04"
$ws.Range("E13").Value = "  +2.81%  "
Set-TextValue $ws.Range("D14") "6.312"
$ws.Range("E14").Value = "  +4.41%  "
Set-TextValue $ws.Range("D15") "7.486"
$ws.Range("E15").Value = "  +4.25%  "
Set-TextValue $ws.Range("D16") "1.817.16"
$ws.Range("E16").Value = "  +3.58%  "
Set-TextValue $ws.Range("D17") "93.52"
$ws.Range("E17").Value = "  +13.34%  "
Set-TextValue $ws.Range("D18") "0.00001084"
$ws.Range("E18").Value = "  +3.37%  "
Set-TextValue $ws.Range("D19") "0.06508"
$ws.Range("E19").Value = "  +9.61%  "
Set-TextValue $ws.Range("D20") "0.9977"
$ws.Range("E20").Value = "  -1.44%  "
Set-TextValue $ws.Range("D21") "17.52"
$ws.Range("E21").Value = "  +5.12%  "
Set-TextValue $ws.Range("D22") "6.251"
$ws.Range("E22").Value = "  +3.41%  "
Set-TextValue $ws.Range("D23") "0.5350"
$ws.Range("E23").Value = "  -1.57%  "
Set-TextValue $ws.Range("D24") "28.303.33"
$ws.Range("E24").Value = "  +3.42%  "
Set-TextValue $ws.Range("D25") "11.73"
$ws.Range("E25").Value = "  +5.31%  "
Set-TextValue $ws.Range("D26") "2.147"
$ws.Range("E26").Value = "  -10.69%  "
Set-TextValue $ws.Range("D27") "20.58"
$ws.Range("E27").Value = "  +3.91%  "
Set-TextValue $ws.Range("D28") "155.44"
$ws.Range("E28").Value = "  +3.79%  "
Set-TextValue $ws.Range("D29") "2.331"
$ws.Range("E29").Value = "  +0.35%  "
Set-TextValue $ws.Range("D30") "2.021.55"
$ws.Range("E30").Value = "  +3.22%  "
Set-TextValue $ws.Range("D31") "127.50"
$ws.Range("E31").Value = "  +1.42%  "
Set-TextValue $ws.Range("D32") "1.202"
$ws.Range("E32").Value = "  -3.50%  "
Set-TextValue $ws.Range("D33") "5.847"
$ws.Range("E33").Value = "  +6.68%  "
Set-TextValue $ws.Range("D34") "0.09238"
$ws.Range("E34").Value = "  +3.06%  "
Set-TextValue $ws.Range("D35") "3.661"
$ws.Range("E35").Value = "  -1.05%  "
Set-TextValue $ws.Range("D36") "13.05"
$ws.Range("E36").Value = "  +6.05%  "
Set-TextValue $ws.Range("D37") "0.02336"
$ws.Range("E37").Value = "  +4.14%  "
Set-TextValue $ws.Range("D38") "0.2172"
$ws.Range("E38").Value = "  +1.73%  "
Set-TextValue $ws.Range("D39") "5.175"
$ws.Range("E39").Value = "  +4.39%  "
Set-TextValue $ws.Range("D40") "0.6559"
$ws.Range("E40").Value = "  +2.97%  "
Set-TextValue $ws.Range("D41") "0.06197"
$ws.Range("E41").Value = "  +2.06%  "
Set-TextValue $ws.Range("D42") "1.193"
$ws.Range("E42").Value = "  +2.34%  "
Set-TextValue $ws.Range("D43") "8.068"
$ws.Range("E43").Value = "  +1.85%  "
Set-TextValue $ws.Range("D44") "0.9970"
$ws.Range("E44").Value = "  -1.58%  "
Set-TextValue $ws.Range("D45") "13.99"
$ws.Range("E45").Value = "  +3.31%  "
Set-TextValue $ws.Range("D46") "1.393"
$ws.Range("E46").Value = "  -2.20%  "
Set-TextValue $ws.Range("D47") "0.6075"
$ws.Range("E47").Value = "  +4.36%  "
Set-TextValue $ws.Range("D48") "3.763"
$ws.Range("E48").Value = "  +0.01%  "
Set-TextValue $ws.Range("D49") "126.77"
$ws.Range("E49").Value = "  +3.20%  "
Set-TextValue $ws.Range("D50") "2.032"
$ws.Range("E50").Value = "  +5.63%  "
Set-TextValue $ws.Range("D51") "0.06990"
$ws.Range("E51").Value = "  +2.09%  "
